# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 1391
$ws.Range("F4").Value  = 13343
$ws.Range("F5").Value  = 766
$ws.Range("F10").Value = 1912
$ws.Range("F13").Value = 20405
$ws.Range("G13").Value = "暂时售罄"
$ws.Range("F14").Value = 541
$ws.Range("F16").Value = 274
$ws.Range("F17").Value = 135
$ws.Range("F20").Value = 314
$ws.Range("F21").Value = 160
$ws.Range("F22").Value = 140
$ws.Range("F24").Value = 234
$ws.Range("F25").Value = 281
$ws.Range("F26").Value = 15
$ws.Range("F27").Value = 1356
$ws.Range("F28").Value = 58

# ---------------------------------------------------------------------
# Sheet: 演出 (Performances)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 4475
$ws.Range("F5").Value  = 190
$ws.Range("F7").Value  = 12
$ws.Range("F9").Value  = 88
$ws.Range("F10").Value = 88
$ws.Range("F11").Value = 385

# ---------------------------------------------------------------------
# Sheet: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 905
$ws.Range("F3").Value = 4426
$ws.Range("F4").Value = 98

# ---------------------------------------------------------------------
# Sheet: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 905
$ws.Range("F5").Value  = 1391
$ws.Range("F6").Value  = 13343
$ws.Range("F8").Value  = 766
$ws.Range("F9").Value  = 4426
$ws.Range("F13").Value = 1912
$ws.Range("F16").Value = 98
$ws.Range("F17").Value = 20406
$ws.Range("G17").Value = "暂时售罄"
$ws.Range("F18").Value = 542
$ws.Range("F19").Value = 4475
$ws.Range("F21").Value = 190
$ws.Range("F22").Value = 190
$ws.Range("F23").Value = 275
$ws.Range("F25").Value = 135
$ws.Range("F26").Value = 12
$ws.Range("F28").Value = 88
$ws.Range("F29").Value = 385
$ws.Range("F32").Value = 314
$ws.Range("F33").Value = 160
$ws.Range("F34").Value = 140
$ws.Range("F37").Value = 234
$ws.Range("F40").Value = 281
$ws.Range("F41").Value = 15
$ws.Range("F42").Value = 1356
$ws.Range("F43").Value = 58
